$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the bill-form header fields (name, designation, year, term, department)
$ws.Range("A3").Value = "নাম: Mr. Md Mehrab Hossain Opi"
$ws.Range("A4").Value = "পদবী: প্রভাষক"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"
$ws.Range("B5").Value = "সিএসই"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# Amount in words for the total bill
$ws.Range("A32").Value = "কথায়:সাত হাজার বাহান্ন টাকা মাত্র।"

# Widen column A so the name/designation text is visible (matches author's resize)
$ws.Columns.Item(1).ColumnWidth = 13.5

# Move the active selection to the totals cell, as recorded by the author
$ws.Range("I32").Select() | Out-Null
